$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.080435
$ws.Range("H2").Value = 24.241305
$ws.Range("I2").Value = 0.1496988574979475
$ws.Range("J2").Value = 0.1496988574979476
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.060105
$ws.Range("N2").Value = 0.180315
$ws.Range("O2").Value = 0.001098708471132188
$ws.Range("P2").Value = 0.001098708471132188
$ws.Range("Q2").Value = 0.485674545675
$ws.Range("R2").Value = 4.371070911075
$ws.Range("S2").Value = 0.0001644754028518052
$ws.Range("T2").Value = 0.0001644754028518053
$ws.Range("G3").Value = 8.080435
$ws.Range("H3").Value = 24.241305
$ws.Range("I3").Value = 0.1496988574979475
$ws.Range("J3").Value = 0.1496988574979476
$ws.Range("O3").Value = 0.002609139380189131
$ws.Range("P3").Value = 0.002609139380189131
$ws.Range("Q3").Value = 1.153347422333333
$ws.Range("R3").Value = 10.380126801
$ws.Range("S3").Value = 0.0003905851842672158
$ws.Range("T3").Value = 0.0003905851842672159
$ws.Range("G4").Value = 8.080435
$ws.Range("H4").Value = 24.241305
$ws.Range("I4").Value = 0.1496988574979475
$ws.Range("J4").Value = 0.1496988574979476
$ws.Range("M4").Value = 33.23770833333333
$ws.Range("N4").Value = 99.71312499999999
$ws.Range("O4").Value = 0.6075792647342859
$ws.Range("P4").Value = 0.607579264734286
$ws.Range("Q4").Value = 268.5751417364583
$ws.Range("R4").Value = 2417.176275628125
$ws.Range("S4").Value = 0.0909539217701656
$ws.Range("T4").Value = 0.09095392177016565
$ws.Range("G5").Value = 8.080435
$ws.Range("H5").Value = 24.241305
$ws.Range("I5").Value = 0.1496988574979475
$ws.Range("J5").Value = 0.1496988574979476
$ws.Range("M5").Value = 0.1028146666666667
$ws.Range("N5").Value = 0.3084440000000001
$ws.Range("O5").Value = 0.001879433411917459
$ws.Range("P5").Value = 0.00187943341191746
$ws.Range("Q5").Value = 0.8307872310466667
$ws.Range("R5").Value = 7.477085079420001
$ws.Range("S5").Value = 0.0002813490345075131
$ws.Range("T5").Value = 0.0002813490345075132
$ws.Range("G6").Value = 8.080435
$ws.Range("H6").Value = 24.241305
$ws.Range("I6").Value = 0.1496988574979475
$ws.Range("J6").Value = 0.1496988574979476
$ws.Range("M6").Value = 21.161778
$ws.Range("N6").Value = 63.485334
$ws.Range("O6").Value = 0.3868334540024753
$ws.Range("P6").Value = 0.3868334540024753
$ws.Range("Q6").Value = 170.99637161343
$ws.Range("R6").Value = 1538.96734452087
$ws.Range("S6").Value = 0.05790852610615538
$ws.Range("T6").Value = 0.05790852610615541
$ws.Range("I7").Value = 0.2404784903431001
$ws.Range("J7").Value = 0.2404784903431001
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.060105
$ws.Range("N7").Value = 0.180315
$ws.Range("O7").Value = 0.001098708471132188
$ws.Range("P7").Value = 0.001098708471132188
$ws.Range("Q7").Value = 0.7801948758600001
$ws.Range("R7").Value = 7.021753882740001
$ws.Range("S7").Value = 0.0002642157544650442
$ws.Range("T7").Value = 0.0002642157544650442
$ws.Range("I8").Value = 0.2404784903431001
$ws.Range("J8").Value = 0.2404784903431001
$ws.Range("O8").Value = 0.002609139380189131
$ws.Range("P8").Value = 0.002609139380189131
$ws.Range("S8").Value = 0.0006274418992426141
$ws.Range("T8").Value = 0.0006274418992426141
$ws.Range("I9").Value = 0.2404784903431001
$ws.Range("J9").Value = 0.2404784903431001
$ws.Range("M9").Value = 33.23770833333333
$ws.Range("N9").Value = 99.71312499999999
$ws.Range("O9").Value = 0.6075792647342859
$ws.Range("P9").Value = 0.607579264734286
$ws.Range("Q9").Value = 431.4431366275001
$ws.Range("R9").Value = 3882.9882296475
$ws.Range("S9").Value = 0.1461097443470719
$ws.Range("T9").Value = 0.1461097443470719
$ws.Range("I10").Value = 0.2404784903431001
$ws.Range("J10").Value = 0.2404784903431001
$ws.Range("M10").Value = 0.1028146666666667
$ws.Range("N10").Value = 0.3084440000000001
$ws.Range("O10").Value = 0.001879433411917459
$ws.Range("P10").Value = 0.00187943341191746
$ws.Range("Q10").Value = 1.334589070736
$ws.Range("R10").Value = 12.011301636624
$ws.Range("S10").Value = 0.0004519633095982925
$ws.Range("T10").Value = 0.0004519633095982926
$ws.Range("I11").Value = 0.2404784903431001
$ws.Range("J11").Value = 0.2404784903431001
$ws.Range("M11").Value = 21.161778
$ws.Range("N11").Value = 63.485334
$ws.Range("O11").Value = 0.3868334540024753
$ws.Range("P11").Value = 0.3868334540024753
$ws.Range("Q11").Value = 274.6911365058961
$ws.Range("R11").Value = 2472.220228553064
$ws.Range("S11").Value = 0.09302512503272233
$ws.Range("T11").Value = 0.09302512503272234
$ws.Range("G12").Value = 15.25749233333333
$ws.Range("H12").Value = 45.772477
$ws.Range("I12").Value = 0.2826616599952471
$ws.Range("J12").Value = 0.2826616599952471
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.060105
$ws.Range("N12").Value = 0.180315
$ws.Range("O12").Value = 0.001098708471132188
$ws.Range("P12").Value = 0.001098708471132188
$ws.Range("Q12").Value = 0.917051576695
$ws.Range("R12").Value = 8.253464190255
$ws.Range("S12").Value = 0.0003105627603010643
$ws.Range("T12").Value = 0.0003105627603010643
$ws.Range("G13").Value = 15.25749233333333
$ws.Range("H13").Value = 45.772477
$ws.Range("I13").Value = 0.2826616599952471
$ws.Range("J13").Value = 0.2826616599952471
$ws.Range("O13").Value = 0.002609139380189131
$ws.Range("P13").Value = 0.002609139380189131
$ws.Range("Q13").Value = 2.177752739044445
$ws.Range("R13").Value = 19.5997746514
$ws.Range("S13").Value = 0.0007375036683632297
$ws.Range("T13").Value = 0.0007375036683632297
$ws.Range("G14").Value = 15.25749233333333
$ws.Range("H14").Value = 45.772477
$ws.Range("I14").Value = 0.2826616599952471
$ws.Range("J14").Value = 0.2826616599952471
$ws.Range("M14").Value = 33.23770833333333
$ws.Range("N14").Value = 99.71312499999999
$ws.Range("O14").Value = 0.6075792647342859
$ws.Range("P14").Value = 0.607579264734286
$ws.Range("Q14").Value = 507.1240800734028
$ws.Range("R14").Value = 4564.116720660625
$ws.Range("S14").Value = 0.1717393635484849
$ws.Range("T14").Value = 0.171739363548485
$ws.Range("G15").Value = 15.25749233333333
$ws.Range("H15").Value = 45.772477
$ws.Range("I15").Value = 0.2826616599952471
$ws.Range("J15").Value = 0.2826616599952471
$ws.Range("M15").Value = 0.1028146666666667
$ws.Range("N15").Value = 0.3084440000000001
$ws.Range("O15").Value = 0.001879433411917459
$ws.Range("P15").Value = 0.00187943341191746
$ws.Range("Q15").Value = 1.568693988420889
$ws.Range("R15").Value = 14.118245895788
$ws.Range("S15").Value = 0.00053124376806312
$ws.Range("T15").Value = 0.0005312437680631201
$ws.Range("G16").Value = 15.25749233333333
$ws.Range("H16").Value = 45.772477
$ws.Range("I16").Value = 0.2826616599952471
$ws.Range("J16").Value = 0.2826616599952471
$ws.Range("M16").Value = 21.161778
$ws.Range("N16").Value = 63.485334
$ws.Range("O16").Value = 0.3868334540024753
$ws.Range("P16").Value = 0.3868334540024753
$ws.Range("Q16").Value = 322.875665594702
$ws.Range("R16").Value = 2905.880990352318
$ws.Range("S16").Value = 0.1093429862500347
$ws.Range("T16").Value = 0.1093429862500347
$ws.Range("G17").Value = 4.142925
$ws.Range("H17").Value = 12.428775
$ws.Range("I17").Value = 0.07675219702895753
$ws.Range("J17").Value = 0.07675219702895753
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.060105
$ws.Range("N17").Value = 0.180315
$ws.Range("O17").Value = 0.001098708471132188
$ws.Range("P17").Value = 0.001098708471132188
$ws.Range("Q17").Value = 0.249010507125
$ws.Range("R17").Value = 2.241094564125
$ws.Range("S17").Value = 0.0000843282890537223872845
$ws.Range("T17").Value = 0.000084328289053722400837
$ws.Range("G18").Value = 4.142925
$ws.Range("H18").Value = 12.428775
$ws.Range("I18").Value = 0.07675219702895753
$ws.Range("J18").Value = 0.07675219702895753
$ws.Range("O18").Value = 0.002609139380189131
$ws.Range("P18").Value = 0.002609139380189131
$ws.Range("Q18").Value = 0.5913334950000001
$ws.Range("R18").Value = 5.322001455000001
$ws.Range("S18").Value = 0.0002002571797842883
$ws.Range("T18").Value = 0.0002002571797842883
$ws.Range("G19").Value = 4.142925
$ws.Range("H19").Value = 12.428775
$ws.Range("I19").Value = 0.07675219702895753
$ws.Range("J19").Value = 0.07675219702895753
$ws.Range("M19").Value = 33.23770833333333
$ws.Range("N19").Value = 99.71312499999999
$ws.Range("O19").Value = 0.6075792647342859
$ws.Range("P19").Value = 0.607579264734286
$ws.Range("Q19").Value = 137.701332796875
$ws.Range("R19").Value = 1239.311995171875
$ws.Range("S19").Value = 0.04663304343759506
$ws.Range("T19").Value = 0.04663304343759506
$ws.Range("G20").Value = 4.142925
$ws.Range("H20").Value = 12.428775
$ws.Range("I20").Value = 0.07675219702895753
$ws.Range("J20").Value = 0.07675219702895753
$ws.Range("M20").Value = 0.1028146666666667
$ws.Range("N20").Value = 0.3084440000000001
$ws.Range("O20").Value = 0.001879433411917459
$ws.Range("P20").Value = 0.00187943341191746
$ws.Range("Q20").Value = 0.4259534529
$ws.Range("R20").Value = 3.833581076100001
$ws.Range("S20").Value = 0.0001442506435342947
$ws.Range("T20").Value = 0.0001442506435342948
$ws.Range("G21").Value = 4.142925
$ws.Range("H21").Value = 12.428775
$ws.Range("I21").Value = 0.07675219702895753
$ws.Range("J21").Value = 0.07675219702895753
$ws.Range("M21").Value = 21.161778
$ws.Range("N21").Value = 63.485334
$ws.Range("O21").Value = 0.3868334540024753
$ws.Range("P21").Value = 0.3868334540024753
$ws.Range("Q21").Value = 87.67165912065001
$ws.Range("R21").Value = 789.0449320858501
$ws.Range("S21").Value = 0.02969031747899016
$ws.Range("T21").Value = 0.02969031747899016
$ws.Range("G22").Value = 13.51654933333334
$ws.Range("H22").Value = 40.549648
$ws.Range("I22").Value = 0.2504087951347477
$ws.Range("J22").Value = 0.2504087951347477
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.060105
$ws.Range("N22").Value = 0.180315
$ws.Range("O22").Value = 0.001098708471132188
$ws.Range("P22").Value = 0.001098708471132188
$ws.Range("Q22").Value = 0.8124121976800001
$ws.Range("R22").Value = 7.311709779120001
$ws.Range("S22").Value = 0.0002751262644605519
$ws.Range("T22").Value = 0.000275126264460552
$ws.Range("G23").Value = 13.51654933333334
$ws.Range("H23").Value = 40.549648
$ws.Range("I23").Value = 0.2504087951347477
$ws.Range("J23").Value = 0.2504087951347477
$ws.Range("O23").Value = 0.002609139380189131
$ws.Range("P23").Value = 0.002609139380189131
$ws.Range("Q23").Value = 1.929262141511112
$ws.Range("R23").Value = 17.3633592736
$ws.Range("S23").Value = 0.0006533514485317825
$ws.Range("T23").Value = 0.0006533514485317825
$ws.Range("G24").Value = 13.51654933333334
$ws.Range("H24").Value = 40.549648
$ws.Range("I24").Value = 0.2504087951347477
$ws.Range("J24").Value = 0.2504087951347477
$ws.Range("M24").Value = 33.23770833333333
$ws.Range("N24").Value = 99.71312499999999
$ws.Range("O24").Value = 0.6075792647342859
$ws.Range("P24").Value = 0.607579264734286
$ws.Range("Q24").Value = 449.2591244144445
$ws.Range("R24").Value = 4043.33211973
$ws.Range("S24").Value = 0.1521431916309684
$ws.Range("T24").Value = 0.1521431916309685
$ws.Range("G25").Value = 13.51654933333334
$ws.Range("H25").Value = 40.549648
$ws.Range("I25").Value = 0.2504087951347477
$ws.Range("J25").Value = 0.2504087951347477
$ws.Range("M25").Value = 0.1028146666666667
$ws.Range("N25").Value = 0.3084440000000001
$ws.Range("O25").Value = 0.001879433411917459
$ws.Range("P25").Value = 0.00187943341191746
$ws.Range("Q25").Value = 1.389699514190223
$ws.Range("R25").Value = 12.507295627712
$ws.Range("S25").Value = 0.0004706266562142389
$ws.Range("T25").Value = 0.000470626656214239
$ws.Range("G26").Value = 13.51654933333334
$ws.Range("H26").Value = 40.549648
$ws.Range("I26").Value = 0.2504087951347477
$ws.Range("J26").Value = 0.2504087951347477
$ws.Range("M26").Value = 21.161778
$ws.Range("N26").Value = 63.485334
$ws.Range("O26").Value = 0.3868334540024753
$ws.Range("P26").Value = 0.3868334540024753
$ws.Range("Q26").Value = 286.0342163180481
$ws.Range("R26").Value = 2574.307946862432
$ws.Range("S26").Value = 0.09686649913457267
$ws.Range("T26").Value = 0.09686649913457268
